# ---------------------------------------------------------------
# Update header text (week number + date range)
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 32   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/12/2025  Through  5/18/2025"

# ---------------------------------------------------------------
# Update numeric statistics cells
# ---------------------------------------------------------------
$numericUpdates = @{
    "D16" = 4
    "E16" = -100
    "F16" = 4
    "G16" = 8
    "H16" = -50
    "J16" = 36
    "K16" = -16.666666666666
    "L16" = -28.571428571428
    "M16" = -30.232558139534
    "N16" = -87.341772151898
    "C17" = 3
    "G17" = 14
    "H17" = 28.571428571428
    "I17" = 57
    "K17" = 11.764705882352
    "L17" = 54.054054054054
    "M17" = 128
    "N17" = -38.043478260869
    "D18" = 4
    "E18" = -100
    "F18" = 7
    "G18" = 15
    "H18" = -53.333333333333
    "I18" = 54
    "J18" = 57
    "K18" = -5.263157894736
    "L18" = -31.645569620253
    "M18" = -25
    "N18" = -87.412587412587
    "C19" = 8
    "D19" = 9
    "E19" = -11.111111111111
    "G19" = 49
    "H19" = -14.285714285714
    "I19" = 192
    "J19" = 228
    "K19" = -15.789473684210
    "L19" = -21.632653061224
    "M19" = 79.439252336448
    "N19" = 38.129496402877
    "G20" = 8
    "H20" = 25
    "L20" = -43.636363636363
    "M20" = -36.734693877551
    "N20" = -90.882352941176
    "C21" = 11
    "D21" = 17
    "E21" = -35.294117647058
    "F21" = 82
    "G21" = 94
    "H21" = -12.765957446808
    "I21" = 365
    "J21" = 404
    "K21" = -9.653465346534
    "L21" = -20.652173913043
    "M21" = 22.895622895622
    "M22" = 66.666666666666
    "F23" = 2
    "H23" = 0
    "M23" = -25
    "C24" = 17
    "D24" = 15
    "E24" = 13.333333333333
    "F24" = 84
    "G24" = 67
    "H24" = 25.373134328358
    "I24" = 373
    "J24" = 360
    "K24" = 3.611111111111
    "L24" = 12.688821752265
    "M24" = 100.537634408602
    "C25" = 9
    "D25" = 9
    "E25" = 0
    "F25" = 48
    "G25" = 36
    "H25" = 33.333333333333
    "I25" = 214
    "J25" = 208
    "K25" = 2.884615384615
    "L25" = 17.582417582417
    "C26" = 3
    "E26" = 50
    "G26" = 13
    "H26" = 69.230769230769
    "I26" = 94
    "J26" = 88
    "K26" = 6.818181818181
    "L26" = 18.987341772151
    "M26" = 23.684210526315
    "I31" = 4
    "L31" = 0
    "L33" = 0
}

foreach ($addr in $numericUpdates.Keys) {
    $ws.Range($addr).Value = $numericUpdates[$addr]
}

# ---------------------------------------------------------------
# Cells that become blank-style text placeholders ("0")
# Copy formatting+value from an existing General-styled "0" cell (C14)
# so the result matches shared string 20 with style 13 (General).
# ---------------------------------------------------------------
$zeroTextCells = @(
    "G15",
    "C16",
    "D17",
    "C18",
    "C20",
    "D20",
    "D23",
    "G27",
    "C28",
    "G28",
)
foreach ($addr in $zeroTextCells) {
    $ws.Range("C14").Copy($ws.Range($addr)) | Out-Null
}

# ---------------------------------------------------------------
# Cells that become "***.*" text placeholders
# Copy formatting+value from an existing General-styled "***.*" cell (E14)
# ---------------------------------------------------------------
$starTextCells = @(
    "H15",
    "E17",
    "E20",
    "E23",
    "H27",
    "H28",
)
foreach ($addr in $starTextCells) {
    $ws.Range("E14").Copy($ws.Range($addr)) | Out-Null
}
